$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2834.5
$ws.Range("I131").Value = 2834.5
$ws.Range("K131").Value = 8503.5
$ws.Range("M131").Value = -3463.5
$ws.Range("H135").Value = 6098053.5
$ws.Range("I135").Value = 440.05713
$ws.Range("K135").Value = 3960.51417
$ws.Range("M135").Value = -1425.51417
$ws.Range("H138").Value = 31255314
$ws.Range("I138").Value = 1162
$ws.Range("J138").Value = 50007810
$ws.Range("K138").Value = 3486
$ws.Range("L138").Value = 150023430
$ws.Range("M138").Value = 1654
$ws.Range("N138").Value = -150033710

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 973.05
$ws.Range("I32").Value = 898.5263
$ws.Range("J32").Value = 2389
$ws.Range("K32").Value = 898.5263
$ws.Range("L32").Value = 2389
$ws.Range("M32").Value = -611.5263
$ws.Range("N32").Value = -2963
$ws.Range("H45").Value = 2802.2068
$ws.Range("I45").Value = 2438.2727
$ws.Range("J45").Value = 3946
$ws.Range("K45").Value = 2438.2727
$ws.Range("L45").Value = 3946
$ws.Range("M45").Value = -2061.2727
$ws.Range("N45").Value = -4700
$ws.Range("H61").Value = 21278892
$ws.Range("I61").Value = 23257792
$ws.Range("J61").Value = 5717.5
$ws.Range("K61").Value = 23257792
$ws.Range("L61").Value = 5717.5
$ws.Range("M61").Value = -23257580
$ws.Range("N61").Value = -6141.5
$ws.Range("H74").Value = 58890572
$ws.Range("I74").Value = 77010130
$ws.Range("K74").Value = 77010130
$ws.Range("M74").Value = -77009256
$ws.Range("H77").Value = 58890572
$ws.Range("I77").Value = 77010130
$ws.Range("K77").Value = 385050650
$ws.Range("M77").Value = -385046282
$ws.Range("H88").Value = 7326.7334
$ws.Range("I88").Value = 8881.75
$ws.Range("J88").Value = 1106.6666
$ws.Range("K88").Value = 8881.75
$ws.Range("L88").Value = 1106.6666
$ws.Range("M88").Value = -8475.75
$ws.Range("N88").Value = -1918.6666
$ws.Range("H91").Value = 7326.7334
$ws.Range("I91").Value = 8881.75
$ws.Range("J91").Value = 1106.6666
$ws.Range("K91").Value = 8881.75
$ws.Range("L91").Value = 1106.6666
$ws.Range("M91").Value = -7477.75
$ws.Range("N91").Value = -3914.6666
$ws.Range("H102").Value = 1583.1305
$ws.Range("J102").Value = 998.3333
$ws.Range("L102").Value = 998.3333
$ws.Range("N102").Value = -4242.3333
$ws.Range("H110").Value = 22680.812
$ws.Range("I110").Value = 27734.385
$ws.Range("J110").Value = 782
$ws.Range("K110").Value = 27734.385
$ws.Range("L110").Value = 782
$ws.Range("M110").Value = -25689.385
$ws.Range("N110").Value = -4872
$ws.Range("H122").Value = 2841.5
$ws.Range("I122").Value = 1887.25
$ws.Range("K122").Value = 5661.75
$ws.Range("M122").Value = -3211.75
$ws.Range("H132").Value = 16435908
$ws.Range("I132").Value = 3440.84
$ws.Range("K132").Value = 10322.52
$ws.Range("M132").Value = -7792.52
$ws.Range("H136").Value = 21278892
$ws.Range("I136").Value = 23257792
$ws.Range("J136").Value = 5717.5
$ws.Range("K136").Value = 69773376
$ws.Range("L136").Value = 17152.5
$ws.Range("M136").Value = -69770826
$ws.Range("N136").Value = -22252.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1519.6774
$ws.Range("I20").Value = 1583.1904
$ws.Range("K20").Value = 1583.1904
$ws.Range("M20").Value = -1336.1904
$ws.Range("H94").Value = 1453.2703
$ws.Range("I94").Value = 1017.4815
$ws.Range("K94").Value = 1017.4815
$ws.Range("M94").Value = -566.4815
$ws.Range("H99").Value = 6029.909
$ws.Range("I99").Value = 4883.8
$ws.Range("J99").Value = 6985
$ws.Range("K99").Value = 4883.8
$ws.Range("L99").Value = 6985
$ws.Range("M99").Value = -3385.8
$ws.Range("N99").Value = -9981
$ws.Range("H105").Value = 7737
$ws.Range("I105").Value = 9985.833000000001
$ws.Range("J105").Value = 2339.8
$ws.Range("K105").Value = 9985.833000000001
$ws.Range("L105").Value = 2339.8
$ws.Range("M105").Value = -8238.833000000001
$ws.Range("N105").Value = -5833.8
$ws.Range("H107").Value = 2812.85
$ws.Range("I107").Value = 2455.4
$ws.Range("K107").Value = 2455.4
$ws.Range("M107").Value = -535.4000000000001
$ws.Range("H134").Value = 2634026.8
$ws.Range("I134").Value = 2779750.5
$ws.Range("J134").Value = 11000
$ws.Range("K134").Value = 8339251.5
$ws.Range("L134").Value = 33000
$ws.Range("M134").Value = -8336716.5
$ws.Range("N134").Value = -38070

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1245.8334
$ws.Range("I16").Value = 1140
$ws.Range("K16").Value = 1140
$ws.Range("M16").Value = -853
$ws.Range("H31").Value = 29414920
$ws.Range("J31").Value = 62503950
$ws.Range("L31").Value = 62503950
$ws.Range("N31").Value = -62504540
$ws.Range("H34").Value = 29414920
$ws.Range("J34").Value = 62503950
$ws.Range("L34").Value = 62503950
$ws.Range("N34").Value = -62504354
$ws.Range("H58").Value = 5535.294
$ws.Range("I58").Value = 5412.5
$ws.Range("K58").Value = 5412.5
$ws.Range("M58").Value = -5209.5
$ws.Range("H107").Value = 619.3077
$ws.Range("I107").Value = 596.4
$ws.Range("K107").Value = 596.4
$ws.Range("M107").Value = 1323.6
$ws.Range("H113").Value = 1245.8334
$ws.Range("I113").Value = 1140
$ws.Range("K113").Value = 1140
$ws.Range("M113").Value = 1030
$ws.Range("H131").Value = 36799.25
$ws.Range("I131").Value = 13000
$ws.Range("K131").Value = 13000
$ws.Range("M131").Value = -7960
$ws.Range("H132").Value = 55064.49
$ws.Range("I132").Value = 79500.58
$ws.Range("K132").Value = 238501.74
$ws.Range("M132").Value = -235971.74
$ws.Range("H136").Value = 5535.294
$ws.Range("I136").Value = 5412.5
$ws.Range("K136").Value = 16237.5
$ws.Range("M136").Value = -13687.5
$ws.Range("H141").Value = 114241.875
$ws.Range("J141").Value = 141489.17
$ws.Range("L141").Value = 141489.17
$ws.Range("N141").Value = -151849.17

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 79.76922999999999
$ws.Range("I38").Value = 38.333332
$ws.Range("J38").Value = 115.28571
$ws.Range("K38").Value = 114.999996
$ws.Range("L38").Value = 345.85713
$ws.Range("M38").Value = 232.000004
$ws.Range("N38").Value = -1039.85713
$ws.Range("H40").Value = 87
$ws.Range("J40").Value = 298.33334
$ws.Range("L40").Value = 1193.33336
$ws.Range("N40").Value = -1331.33336
$ws.Range("H107").Value = 645.5278
$ws.Range("I107").Value = 468.46155
$ws.Range("J107").Value = 745.6087
$ws.Range("K107").Value = 1405.38465
$ws.Range("L107").Value = 2236.8261
$ws.Range("M107").Value = 514.61535
$ws.Range("N107").Value = -6076.8261
$ws.Range("H131").Value = 1479.1177
$ws.Range("J131").Value = 1708.375
$ws.Range("L131").Value = 5125.125
$ws.Range("N131").Value = -15205.125

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 129999.5
$ws.Range("J108").Value = 129999.5
$ws.Range("L108").Value = 129999.5
$ws.Range("N108").Value = -137679.5
$ws.Range("H122").Value = 2546.3333
$ws.Range("I122").Value = 2599.8
$ws.Range("J122").Value = 2439.4
$ws.Range("K122").Value = 7799.400000000001
$ws.Range("L122").Value = 7318.200000000001
$ws.Range("M122").Value = -5349.400000000001
$ws.Range("N122").Value = -12218.2
$ws.Range("H126").Value = 21503354
$ws.Range("I126").Value = 11225385
$ws.Range("J126").Value = 40003696
$ws.Range("K126").Value = 33676155
$ws.Range("L126").Value = 120011088
$ws.Range("M126").Value = -33673685
$ws.Range("N126").Value = -120016028
$ws.Range("H132").Value = 2507.8975
$ws.Range("I132").Value = 2304.7354
$ws.Range("K132").Value = 6914.206200000001
$ws.Range("M132").Value = -4384.206200000001
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4301.7144
$ws.Range("I40").Value = 4318.6665
$ws.Range("K40").Value = 4318.6665
$ws.Range("M40").Value = -4182.6665
$ws.Range("H46").Value = 1288.8889
$ws.Range("J46").Value = 2987.5
$ws.Range("L46").Value = 2987.5
$ws.Range("N46").Value = -3363.5
$ws.Range("H55").Value = 698.2857
$ws.Range("I55").Value = 249.75
$ws.Range("J55").Value = 877.7
$ws.Range("K55").Value = 249.75
$ws.Range("L55").Value = 877.7
$ws.Range("M55").Value = -76.75
$ws.Range("N55").Value = -1223.7
$ws.Range("H122").Value = 3834.476
$ws.Range("I122").Value = 3164.3076
$ws.Range("J122").Value = 4923.5
$ws.Range("K122").Value = 9492.9228
$ws.Range("L122").Value = 14770.5
$ws.Range("M122").Value = -7042.9228
$ws.Range("N122").Value = -19670.5
$ws.Range("H132").Value = 29417884
$ws.Range("I132").Value = 2943.611
$ws.Range("J132").Value = 142875520
$ws.Range("K132").Value = 8830.832999999999
$ws.Range("L132").Value = 428626560
$ws.Range("M132").Value = -6300.832999999999
$ws.Range("N132").Value = -428631620

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 100001890
$ws.Range("I122").Value = 125001176
$ws.Range("K122").Value = 375003528
$ws.Range("M122").Value = -375001078
$ws.Range("H126").Value = 9362.666999999999
$ws.Range("I126").Value = 10595
$ws.Range("J126").Value = 4433.3335
$ws.Range("K126").Value = 31785
$ws.Range("L126").Value = 13300.0005
$ws.Range("M126").Value = -29315
$ws.Range("N126").Value = -18240.0005
$ws.Range("H136").Value = 2000.5555
$ws.Range("I136").Value = 1041.3469
$ws.Range("J136").Value = 11400.8
$ws.Range("K136").Value = 3124.0407
$ws.Range("L136").Value = 34202.39999999999
$ws.Range("M136").Value = -574.0407
$ws.Range("N136").Value = -39302.39999999999
